$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (date serial, nuovi pos., somma mobile 7gg., somma mobile per 100mila abitanti)
$newRows = @(
    @(44326, 0, 5, 119.1327138432213),
    @(44327, 0, 5, 119.1327138432213),
    @(44328, 0, 5, 119.1327138432213),
    @(44329, 0, 5, 119.1327138432213)
)

$startRow = 252
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Carry over the formatting (date number format / alignment / border / font)
    # from the last existing data row (251) into column A of the new row.
    $ws.Range("A" + ($r - 1)).Copy() | Out-Null
    $ws.Range("A" + $r).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = $false
